# Generate Report for Archive
#
# 1) The "Status" value that was previously "Ready for handoff" is now
#    "In Translation" everywhere it appears (Overview!E2, Overview!F2,
#    zh-cn!C2, de-de!C2 all shared the same string).
# 2) The Status columns got narrower to match the new (shorter) text:
#    Overview columns E & F, and column C on both the "zh-cn" and
#    "de-de" sheets shrink from ~17.216 to ~13.410 (stored OOXML width
#    units). Excel's ColumnWidth property is specified in "characters"
#    and gets stored as characters + 5/6, so we back that offset out of
#    the desired stored width before assigning it.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$targetStoredWidth = 13.4101845877511
$colWidth = $targetStoredWidth - (5/6)

# --- Overview sheet: columns E (zh-cn) and F (de-de), row 2 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $colWidth
$wsOverview.Columns.Item(6).ColumnWidth = $colWidth

# --- zh-cn sheet: Status column C, row 2 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $colWidth

# --- de-de sheet: Status column C, row 2 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $colWidth
